$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the points-table headers (row 1) - tidy up casing / wording
# for the new "Stats" section.
$ws.Range("B1").Value = "Pld"
$ws.Range("C1").Value = "Won"
$ws.Range("D1").Value = "Lost"
$ws.Range("E1").Value = "Tied"
$ws.Range("F1").Value = "Net RR"
$ws.Range("G1").Value = "Pts"

# Select the whole table, matching the workbook's last saved selection.
$ws.Range("A1:G9").Select()
